$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.057090308610182
$ws.Cells.Item(2, 4).Value = 1.062861740563521
$ws.Cells.Item(2, 5).Value = 1.05326764887379
$ws.Cells.Item(2, 6).Value = 1.070573051813128
$ws.Cells.Item(2, 9).Value = 1.041643376295184
$ws.Cells.Item(2, 10).Value = 1.06208854356696
$ws.Cells.Item(2, 11).Value = 1.065581940082656
$ws.Cells.Item(2, 12).Value = 1.056014088241574
$ws.Cells.Item(2, 13).Value = 1.073272529633108
$ws.Cells.Item(2, 14).Value = 1.063596831122084

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.058531525792873
$ws.Cells.Item(3, 4).Value = 1.064171462249029
$ws.Cells.Item(3, 5).Value = 1.054510454566504
$ws.Cells.Item(3, 6).Value = 1.071958476059719
$ws.Cells.Item(3, 9).Value = 1.04193817245341
$ws.Cells.Item(3, 10).Value = 1.063179726572321
$ws.Cells.Item(3, 11).Value = 1.06670531745953
$ws.Cells.Item(3, 12).Value = 1.057068803768172
$ws.Cells.Item(3, 13).Value = 1.074472934449262
$ws.Cells.Item(3, 14).Value = 1.064689563732473

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.059463024570496
$ws.Cells.Item(4, 4).Value = 1.065018175714777
$ws.Cells.Item(4, 5).Value = 1.055313936493429
$ws.Cells.Item(4, 6).Value = 1.072854244321519
$ws.Cells.Item(4, 9).Value = 1.042127031040905
$ws.Cells.Item(4, 10).Value = 1.063884289154497
$ws.Cells.Item(4, 11).Value = 1.067430905465773
$ws.Cells.Item(4, 12).Value = 1.057750019816613
$ws.Cells.Item(4, 13).Value = 1.075248447199305
$ws.Cells.Item(4, 14).Value = 1.06539512687433

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.059854377441247
$ws.Cells.Item(5, 4).Value = 1.065373955635806
$ws.Cells.Item(5, 5).Value = 1.055651557686768
$ws.Cells.Item(5, 6).Value = 1.073230663945246
$ws.Cells.Item(5, 9).Value = 1.042205975000232
$ws.Cells.Item(5, 10).Value = 1.064180130882346
$ws.Cells.Item(5, 11).Value = 1.067735632612595
$ws.Cells.Item(5, 12).Value = 1.058036106296489
$ws.Cells.Item(5, 13).Value = 1.075574182743735
$ws.Cells.Item(5, 14).Value = 1.065691388731366

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.059920072863525
$ws.Cells.Item(6, 4).Value = 1.065433682317604
$ws.Cells.Item(6, 5).Value = 1.05570823635478
$ws.Cells.Item(6, 6).Value = 1.07329385713245
$ws.Cells.Item(6, 9).Value = 1.042219203559124
$ws.Cells.Item(6, 10).Value = 1.064229783230785
$ws.Cells.Item(6, 11).Value = 1.067786779578932
$ws.Cells.Item(6, 12).Value = 1.05808412418411
$ws.Cells.Item(6, 13).Value = 1.075628858293247
$ws.Cells.Item(6, 14).Value = 1.065741111591835

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.059468254818797
$ws.Cells.Item(7, 4).Value = 1.065022930358998
$ws.Cells.Item(7, 5).Value = 1.055318448438171
$ws.Cells.Item(7, 6).Value = 1.072859274687065
$ws.Cells.Item(7, 9).Value = 1.042128087669158
$ws.Cells.Item(7, 10).Value = 1.063888243601403
$ws.Cells.Item(7, 11).Value = 1.067434978458028
$ws.Cells.Item(7, 12).Value = 1.057753843679467
$ws.Cells.Item(7, 13).Value = 1.075252800827819
$ws.Cells.Item(7, 14).Value = 1.065399086937004

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.057577598725911
$ws.Cells.Item(8, 4).Value = 1.06330452805065
$ws.Cells.Item(8, 5).Value = 1.053687807572363
$ws.Cells.Item(8, 6).Value = 1.071041408964861
$ws.Cells.Item(8, 9).Value = 1.041743397033284
$ws.Cells.Item(8, 10).Value = 1.062457628233925
$ws.Cells.Item(8, 11).Value = 1.065961864922033
$ws.Cells.Item(8, 12).Value = 1.05637079653136
$ws.Cells.Item(8, 13).Value = 1.073678469226258
$ws.Cells.Item(8, 14).Value = 1.06396643993161

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.054237604341073
$ws.Cells.Item(9, 4).Value = 1.060270429924645
$ws.Cells.Item(9, 5).Value = 1.05080889042359
$ws.Cells.Item(9, 6).Value = 1.067832569154875
$ws.Cells.Item(9, 9).Value = 1.041050956977599
$ws.Cells.Item(9, 10).Value = 1.059924983934452
$ws.Cells.Item(9, 11).Value = 1.063355815353896
$ws.Cells.Item(9, 12).Value = 1.053923904581825
$ws.Cells.Item(9, 13).Value = 1.070894685820403
$ws.Cells.Item(9, 14).Value = 1.061430198986734

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.052004917945008
$ws.Cells.Item(10, 4).Value = 1.058243342871945
$ws.Cells.Item(10, 5).Value = 1.048885639525954
$ws.Cells.Item(10, 6).Value = 1.065689324313388
$ws.Cells.Item(10, 9).Value = 1.040579453998969
$ws.Cells.Item(10, 10).Value = 1.058228415995733
$ws.Cells.Item(10, 11).Value = 1.06161130672505
$ws.Cells.Item(10, 12).Value = 1.052285822262252
$ws.Cells.Item(10, 13).Value = 1.06903209554026
$ws.Cells.Item(10, 14).Value = 1.059731221726943

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.051036629709987
$ws.Cells.Item(11, 4).Value = 1.057364494790368
$ws.Cells.Item(11, 5).Value = 1.04805184942564
$ws.Cells.Item(11, 6).Value = 1.064760254756787
$ws.Cells.Item(11, 9).Value = 1.040372926629601
$ws.Cells.Item(11, 10).Value = 1.057491794429183
$ws.Cells.Item(11, 11).Value = 1.060854162535515
$ws.Cells.Item(11, 12).Value = 1.05157484233472
$ws.Cells.Item(11, 13).Value = 1.068223914472994
$ws.Cells.Item(11, 14).Value = 1.058993554073277

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.050676728819026
$ws.Cells.Item(12, 4).Value = 1.057037879708207
$ws.Cells.Item(12, 5).Value = 1.047741985787024
$ws.Cells.Item(12, 6).Value = 1.064414996220993
$ws.Cells.Item(12, 9).Value = 1.040295856396289
$ws.Cells.Item(12, 10).Value = 1.057217875297328
$ws.Cells.Item(12, 11).Value = 1.060572655843239
$ws.Cells.Item(12, 12).Value = 1.051310495643311
$ws.Cells.Item(12, 13).Value = 1.067923464170971
$ws.Cells.Item(12, 14).Value = 1.058719245944836

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.050753939554776
$ws.Cells.Item(13, 4).Value = 1.057107947615643
$ws.Cells.Item(13, 5).Value = 1.047808459771406
$ws.Cells.Item(13, 6).Value = 1.064489062718199
$ws.Cells.Item(13, 9).Value = 1.040312404394919
$ws.Cells.Item(13, 10).Value = 1.057276645769538
$ws.Cells.Item(13, 11).Value = 1.060633052268521
$ws.Cells.Item(13, 12).Value = 1.051367210616736
$ws.Cells.Item(13, 13).Value = 1.067987923446282
$ws.Cells.Item(13, 14).Value = 1.058778099877857

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.051006885036108
$ws.Cells.Item(14, 4).Value = 1.057337500222503
$ws.Cells.Item(14, 5).Value = 1.048026239230729
$ws.Cells.Item(14, 6).Value = 1.064731718891918
$ws.Cells.Item(14, 9).Value = 1.040366563267941
$ws.Cells.Item(14, 10).Value = 1.05746915844177
$ws.Cells.Item(14, 11).Value = 1.060830898634723
$ws.Cells.Item(14, 12).Value = 1.051552996627001
$ws.Cells.Item(14, 13).Value = 1.068199084408292
$ws.Cells.Item(14, 14).Value = 1.058970885940165

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.051162701671954
$ws.Cells.Item(15, 4).Value = 1.057478912229639
$ws.Cells.Item(15, 5).Value = 1.048160399401231
$ws.Cells.Item(15, 6).Value = 1.064881205892096
$ws.Cells.Item(15, 9).Value = 1.040399885010505
$ws.Cells.Item(15, 10).Value = 1.057587731271919
$ws.Cells.Item(15, 11).Value = 1.060952762427884
$ws.Cells.Item(15, 12).Value = 1.051667431309524
$ws.Cells.Item(15, 13).Value = 1.068329153579021
$ws.Cells.Item(15, 14).Value = 1.059089627157333

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.05206914740861
$ws.Cells.Item(16, 4).Value = 1.058301645344497
$ws.Cells.Item(16, 5).Value = 1.048940953714216
$ws.Cells.Item(16, 6).Value = 1.065750961385504
$ws.Cells.Item(16, 9).Value = 1.040593110600924
$ws.Cells.Item(16, 10).Value = 1.058277260642248
$ws.Cells.Item(16, 11).Value = 1.061661518362744
$ws.Cells.Item(16, 12).Value = 1.052332971867687
$ws.Cells.Item(16, 13).Value = 1.069085696270805
$ws.Cells.Item(16, 14).Value = 1.059780135738459

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.052637324832394
$ws.Cells.Item(17, 4).Value = 1.058817424032764
$ws.Cells.Item(17, 5).Value = 1.049430300833589
$ws.Cells.Item(17, 6).Value = 1.066296255986814
$ws.Cells.Item(17, 9).Value = 1.040713681953296
$ws.Cells.Item(17, 10).Value = 1.058709246074168
$ws.Cells.Item(17, 11).Value = 1.062105627295932
$ws.Cells.Item(17, 12).Value = 1.052749994877568
$ws.Cells.Item(17, 13).Value = 1.06955980536431
$ws.Cells.Item(17, 14).Value = 1.060212734639242

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.052968586908785
$ws.Cells.Item(18, 4).Value = 1.059118162682963
$ws.Cells.Item(18, 5).Value = 1.049715631652577
$ws.Cells.Item(18, 6).Value = 1.066614218087229
$ws.Cells.Item(18, 9).Value = 1.040783781257923
$ws.Cells.Item(18, 10).Value = 1.0589610235889
$ws.Cells.Item(18, 11).Value = 1.062364499015788
$ws.Cells.Item(18, 12).Value = 1.052993075623433
$ws.Cells.Item(18, 13).Value = 1.069836184685141
$ws.Cells.Item(18, 14).Value = 1.060464869706923

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.053081514132255
$ws.Cells.Item(19, 4).Value = 1.059220688923812
$ws.Cells.Item(19, 5).Value = 1.049812905777699
$ws.Cells.Item(19, 6).Value = 1.066722618330384
$ws.Cells.Item(19, 9).Value = 1.040807644699909
$ws.Cells.Item(19, 10).Value = 1.05904684075567
$ws.Cells.Item(19, 11).Value = 1.062452738985767
$ws.Cells.Item(19, 12).Value = 1.053075932588867
$ws.Cells.Item(19, 13).Value = 1.069930395831292
$ws.Cells.Item(19, 14).Value = 1.060550808743914

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.052576379946071
$ws.Cells.Item(20, 4).Value = 1.058762096900006
$ws.Cells.Item(20, 5).Value = 1.049377808597445
$ws.Cells.Item(20, 6).Value = 1.066237761321375
$ws.Cells.Item(20, 9).Value = 1.040700769375926
$ws.Cells.Item(20, 10).Value = 1.058662918051261
$ws.Cells.Item(20, 11).Value = 1.062057996180687
$ws.Cells.Item(20, 12).Value = 1.052705269017414
$ws.Cells.Item(20, 13).Value = 1.069508954585146
$ws.Cells.Item(20, 14).Value = 1.060166340825228

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.050932405455476
$ws.Cells.Item(21, 4).Value = 1.057269907488263
$ws.Cells.Item(21, 5).Value = 1.047962112980026
$ws.Cells.Item(21, 6).Value = 1.064660267181719
$ws.Cells.Item(21, 9).Value = 1.040350624692007
$ws.Cells.Item(21, 10).Value = 1.05741247672131
$ws.Cells.Item(21, 11).Value = 1.060772645295603
$ws.Cells.Item(21, 12).Value = 1.051498294412969
$ws.Cells.Item(21, 13).Value = 1.068136909874776
$ws.Cells.Item(21, 14).Value = 1.05891412372516

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.049897407195317
$ws.Cells.Item(22, 4).Value = 1.056330712143423
$ws.Cells.Item(22, 5).Value = 1.047071098765994
$ws.Cells.Item(22, 6).Value = 1.063667501072309
$ws.Cells.Item(22, 9).Value = 1.040128409809463
$ws.Cells.Item(22, 10).Value = 1.056624506934358
$ws.Cells.Item(22, 11).Value = 1.059962931619202
$ws.Cells.Item(22, 12).Value = 1.050737931745628
$ws.Cells.Item(22, 13).Value = 1.067272767829404
$ws.Cells.Item(22, 14).Value = 1.05812503493073

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.050446211248638
$ws.Cells.Item(23, 4).Value = 1.056828693717137
$ws.Cells.Item(23, 5).Value = 1.047543530267753
$ws.Cells.Item(23, 6).Value = 1.064193875591391
$ws.Cells.Item(23, 9).Value = 1.040246406430199
$ws.Cells.Item(23, 10).Value = 1.057042393908485
$ws.Cells.Item(23, 11).Value = 1.060392325997131
$ws.Cells.Item(23, 12).Value = 1.05114115720049
$ws.Cells.Item(23, 13).Value = 1.067731008266172
$ws.Cells.Item(23, 14).Value = 1.058543515352292

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.052603918768945
$ws.Cells.Item(24, 4).Value = 1.058787097178689
$ws.Cells.Item(24, 5).Value = 1.049401527879768
$ws.Cells.Item(24, 6).Value = 1.066264192849137
$ws.Cells.Item(24, 9).Value = 1.040706604718223
$ws.Cells.Item(24, 10).Value = 1.058683852285078
$ws.Cells.Item(24, 11).Value = 1.062079519155486
$ws.Cells.Item(24, 12).Value = 1.052725479208419
$ws.Cells.Item(24, 13).Value = 1.069531932359646
$ws.Cells.Item(24, 14).Value = 1.060187304788058

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.055102107777775
$ws.Cells.Item(25, 4).Value = 1.061055563619789
$ws.Cells.Item(25, 5).Value = 1.051553840453558
$ws.Cells.Item(25, 6).Value = 1.068662817196399
$ws.Cells.Item(25, 9).Value = 1.04123170397125
$ws.Cells.Item(25, 10).Value = 1.060581148449073
$ws.Cells.Item(25, 11).Value = 1.064030781295316
$ws.Cells.Item(25, 12).Value = 1.05455766926246
$ws.Cells.Item(25, 13).Value = 1.071615527434336
$ws.Cells.Item(25, 14).Value = 1.062087295330229
